$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated 2D training schedule values for rows 2-5 (columns B:I), and a new row 6.
$data = @(
    @(2, 9, 7, 6, 4, -3, -3, 34, 5),
    @(3, 5, 5, 4, 0, -1, -5, 56, 5),
    @(4, 8, 6, 3, 5, -5, -1, 12, 5),
    @(5, 5, 7, 3, 3, -2, -4, 45, 5),
    @(6, 9, 6, 5, 4, -4, -2, 23, 5)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 1).Value = $row - 1
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $ws.Cells.Item($row, 7).Value = $entry[6]
    $ws.Cells.Item($row, 8).Value = $entry[7]
    $ws.Cells.Item($row, 9).Value = $entry[8]
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

[void]$ws.Range("I1").Select()
